# Apply the "atividade" name updates (infinitive verb forms) on the
# "Mapa geral do Processo" sheet, resize column B to fit the longer
# labels, and shrink the process-map picture's width to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapa geral do Processo")

$ws.Range("B2").Value = "1.1 Elicitar"
$ws.Range("B3").Value = "1.2 Documentar"
$ws.Range("B4").Value = "1.3 Validar Requisito"
$ws.Range("B5").Value = "2.1 Criar Design De Arquitetura "
$ws.Range("B6").Value = "2.2 Projetar BD"
$ws.Range("B7").Value = "3.1 Codificar"
$ws.Range("B8").Value = "3.2 Testar"
$ws.Range("B9").Value = "3.3 Entregar"

# Widen column B so the new (longer) activity labels fit.
$ws.Range("B1").ColumnWidth = 23.8333

# Shrink the process-map picture so its right edge moves in by ~17pt,
# matching the updated column layout. Height / top-left stay put.
# (Width is tuned empirically so the resulting anchor offset lines up;
# it must be set after the column resize above since that shifts the
# pixel grid the picture's right-edge column/offset is measured against.)
$shp = $ws.Shapes.Item(1)
$shp.Width = 549.0874
